$p = $ppt.ActivePresentation

# --- 1) Slide 2 ("Bevezető", id 257): the paragraph after the manual line
#        break (<a:br/>) was split over several runs, one of which had a
#        typo ("maximalizálásámak"). Collapse them into a single run with
#        the corrected text. We only touch the characters that follow the
#        line break (TextRange.Characters(start, length)) so the
#        <a:br/> itself, and everything before it, stay untouched.
$s2 = $p.Slides.Item(2)
$body2 = $s2.Shapes.Item(2)
$tr2 = $body2.TextFrame.TextRange
$oldTail = "Mivel a project még kezdetleges állapotban van, így kollégáink visszajelzései alapján további fejlesztéseket fogunk végrehajtani, a maximális felhasználói élmény maximalizálásámak, valamint a munkafolyamatok megkönnyítésének érdekében."
$newTail = "Mivel a project még kezdetleges állapotban van, így kollégáink visszajelzései alapján további fejlesztéseket fogunk végrehajtani, a maximális felhasználói élmény maximalizálásának, valamint a munkafolyamatok megkönnyítésének érdekében."
$tailStart = $tr2.Length - $oldTail.Length + 1
$tailRange = $tr2.Characters($tailStart, $oldTail.Length)
$tailRange.Text = $newTail

# --- 2) Insert the new "Jelenlegi helyzet" slide as slide 3 (id 258),
#        right after slide 2, reusing the same "Title and Content"
#        layout ( slideLayout2.xml ) that slide 2 itself uses.
$s3 = $p.Slides.Add(3, 2)

$title3 = $s3.Shapes.Item(1)
$title3.Name = "Cím 1"
$title3.TextFrame.TextRange.Text = "Jelenlegi helyzet"
$title3.TextFrame.TextRange.LanguageID = "hu-HU"

$body3 = $s3.Shapes.Item(2)
$body3.Name = "Tartalom helye 2"
$body3.TextFrame.TextRange.Text = "Cégünk jelenleg is használt adminisztrációs rendszerét korábban az Önök cége készítette el számunkra. A programmal elégedettek vagyunk, azonban a gyorsan fejlődő világban, rohamosan változó piaci helyzet mellett, elkerülhetetlenné vált ezen program továbbfejlesztése is. Programunkban képesek vagyunk rögzíteni az autókat valamint a hozzájuk kapcsolódó ügyfeleket. Itt követjük nyomon a szerelések árát is. A kifizetett autókat egy hónapon belül töröljük az adatbázisból. Azt is tudni érdemes, hogy az Önök által használt alkalmazásának is vannak hiányosságai amik orvosolása nagyban megkönnyítené a jelenlegi adminisztrációs folyamatainkat, ezért megkérjük önöket az alkalmazás továbbfejlesztésére."
$body3.TextFrame.TextRange.LanguageID = "hu-HU"
$body3.TextFrame.AutoSize = 2
